$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -1.221066362194925
$ws.Range("D2").Value = 0.2349861455002353

$ws.Range("C3").Value = 0.1659878138695713
$ws.Range("D3").Value = 0.8696825269947404

$ws.Range("C4").Value = -0.4076334686602358
$ws.Range("D4").Value = 0.6874823318318857

$ws.Range("C5").Value = -0.04005864555531555
$ws.Range("D5").Value = 0.9684076859375383

$ws.Range("C6").Value = 1.190858081340474
$ws.Range("D6").Value = 0.2464000520199752

$ws.Range("C7").Value = 1.088535489911793
$ws.Range("D7").Value = 0.2881401356405855

$ws.Range("C8").Value = 1.568383225484835
$ws.Range("D8").Value = 0.1310638319071515

$ws.Range("C9").Value = -0.4900606732739369
$ws.Range("D9").Value = 0.628942557726103

$ws.Range("C10").Value = -0.1627367419473806
$ws.Range("D10").Value = 0.8722112651880236

$ws.Range("C11").Value = 0.3253140268527125
$ws.Range("D11").Value = 0.7480170426688049
